$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Content edit -----------------------------------------------------
# The test method used to validate the "display" button covered only the
# Visits screen. It now also covers the Offices-Visits screen, so the
# TestMethodName cell (A2) records both method names, one per line.
$ws.Range("A2").Value = "ValidateOnVisitsDisplayBtnFunctionality,`nValidateOnOfficesVisitsDisplayBtnFunctionality"

# Wrapped, multi-line text needs a taller row to display fully.
$ws.Rows.Item(2).RowHeight = 90

# --- View-state edit ----------------------------------------------------
# Reviewer scrolled back to the left edge of the sheet (clearing the old
# "topLeftCell=I1" scroll position) and left the selection on B4:B6.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B4:B6").Select() | Out-Null
